$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Alexander, Bradley and Gonzales"

$ws.Range("A2").Value = 'Nihar'
$ws.Range("B2").Value = '2023PCP5317'
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = '8050106439'
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = 'niharkajla28@gmail.com'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1995-10-28'
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = 'Male'
$ws.Range("G2").Value = 'PG'
$ws.Range("H2").Value = 'CSE'
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = '88'
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = '85'
$ws.Range("J2").Style = "Normal"
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = '8.167'
$ws.Range("K2").Style = "Normal"
$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = '1'
$ws.Range("L2").Style = "Normal"
$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = '3'
$ws.Range("M2").Style = "Normal"
$ws.Range("N2").Value = 'General'

$ws.Range("A3").Value = 'Samay Raina'
$ws.Range("B3").Value = '2023PCP5320'
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = '8050106439'
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = 'niharkajla28@gmail.com'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2023-12-07'
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = 'Male'
$ws.Range("G3").Value = 'PG'
$ws.Range("H3").Value = 'PS'
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = '90'
$ws.Range("I3").Style = "Normal"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = '95'
$ws.Range("J3").Style = "Normal"
$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = '8'
$ws.Range("K3").Style = "Normal"
$ws.Range("L3").NumberFormat = "@"
$ws.Range("L3").Value = '1'
$ws.Range("L3").Style = "Normal"
$ws.Range("M3").NumberFormat = "@"
$ws.Range("M3").Value = '0'
$ws.Range("M3").Style = "Normal"
$ws.Range("N3").Value = 'General'

$ws.Range("A4").Value = 'Suhana Sharma'
$ws.Range("B4").Value = '2023PCP5305'
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = '8050106439'
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = 'niharamazon5005@gmail.com'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '2024-03-07'
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = 'Female'
$ws.Range("G4").Value = 'PG'
$ws.Range("H4").Value = 'CSE'
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = '90'
$ws.Range("I4").Style = "Normal"
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = '90'
$ws.Range("J4").Style = "Normal"
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = '9'
$ws.Range("K4").Style = "Normal"
$ws.Range("L4").NumberFormat = "@"
$ws.Range("L4").Value = '0'
$ws.Range("L4").Style = "Normal"
$ws.Range("M4").NumberFormat = "@"
$ws.Range("M4").Value = '0'
$ws.Range("M4").Style = "Normal"
$ws.Range("N4").Value = 'General'

$ws.Range("A5").Value = 'Sagar Shah'
$ws.Range("B5").Value = '2023PCP5319'
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = '8050106439'
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = 'niharkajla123@gmail.com'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2024-02-01'
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = 'Male'
$ws.Range("G5").Value = 'PG'
$ws.Range("H5").Value = 'VLSI'
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = '80'
$ws.Range("I5").Style = "Normal"
$ws.Range("J5").NumberFormat = "@"
$ws.Range("J5").Value = '85'
$ws.Range("J5").Style = "Normal"
$ws.Range("K5").NumberFormat = "@"
$ws.Range("K5").Value = '8.75'
$ws.Range("K5").Style = "Normal"
$ws.Range("L5").NumberFormat = "@"
$ws.Range("L5").Value = '0'
$ws.Range("L5").Style = "Normal"
$ws.Range("M5").NumberFormat = "@"
$ws.Range("M5").Value = '0'
$ws.Range("M5").Style = "Normal"
$ws.Range("N5").Value = 'OBC'

$ws.Range("A6").Value = 'Ellen Degenerous'
$ws.Range("B6").Value = '2023PCP5318'
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = '8050106439'
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = 'niharkajla28@gmail.com'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '2024-03-05'
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = 'Female'
$ws.Range("G6").Value = 'PG'
$ws.Range("H6").Value = 'CSIS'
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = '70'
$ws.Range("I6").Style = "Normal"
$ws.Range("J6").NumberFormat = "@"
$ws.Range("J6").Value = '70'
$ws.Range("J6").Style = "Normal"
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = '9'
$ws.Range("K6").Style = "Normal"
$ws.Range("L6").NumberFormat = "@"
$ws.Range("L6").Value = '0'
$ws.Range("L6").Style = "Normal"
$ws.Range("M6").NumberFormat = "@"
$ws.Range("M6").Value = '0'
$ws.Range("M6").Style = "Normal"
$ws.Range("N6").Value = 'General'

$ws.Columns.Item(1).ColumnWidth = 15.0
